# ---------------------------------------------------------------------------
# api_function_complete_tracker.xlsx edit:
#   1. Sheet1!P16 changes from "Guatemala" to "Guatemala**"
#   2. Sheet1 selection changes to A1:A13
#   3. A brand new sheet "names for taxa" is inserted right after Sheet1,
#      containing a small lookup table of taxonomic-rank field names per
#      data provider / atlas.
# ---------------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Sheet1 edits
# ---------------------------------------------------------------------------
$ws1.Range("P16").Value = "Guatemala**"

# ---------------------------------------------------------------------------
# 2) Insert the new worksheet right after Sheet1
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "names for taxa"

# Stamp every cell in the A1:I13 block with the same (bold, bordered,
# no-fill) style used by Sheet1's header row, by copying that row's
# formatting onto each target row. This reproduces style index 2 for
# every cell in the block, matching the un-highlighted label/data cells
# of the new sheet.
for ($r = 1; $r -le 13; $r++) {
    $ws1.Range("A1:I1").Copy($ws2.Range("A$r"))
}

# Column A (the country/provider label column) re-uses the same three
# label styles already present on Sheet1's own country column (rows
# 2-13): bold+italic (9), bold (7), plain (2).
$ws1.Range("A2:A13").Copy($ws2.Range("A2"))
# ... except GBIF's row, which on the new sheet uses the plain bold
# style (7) rather than Sheet1's special "API provider" style (11).
$ws1.Range("A7").Copy($ws2.Range("A8"))

# --- Row 1: column headers -------------------------------------------------
$ws2.Range("B1").Value = "Kingdom"
$ws2.Range("C1").Value = "Phylum"
$ws2.Range("D1").Value = "Class"
$ws2.Range("E1").Value = "Order"
$ws2.Range("F1").Value = "Family"
$ws2.Range("G1").Value = "Genus"
$ws2.Range("H1").Value = "Species"
$ws2.Range("I1").Value = "Subspecies"

# --- Row 2: Australia --------------------------------------------------
$ws2.Range("A2").Value = "Australia"
$ws2.Range("B2").Value = "kingdom, kingdomID"
$ws2.Range("C2").Value = "phylum, phylumID"
$ws2.Range("D2").Value = "class, classID"
$ws2.Range("E2").Value = "order"
$ws2.Range("F2").Value = "family"
$ws2.Range("G2").Value = "genus"
$ws2.Range("H2").Value = "species, speciesID"
$ws2.Range("I2").Value = "subspecies, subspeciesID"

# --- Row 3: Austria ------------------------------------------------------
$ws2.Range("A3").Value = "Austria"
$ws2.Range("B3").Value = "kingdom"
$ws2.Range("C3").Value = "phylum"
$ws2.Range("D3").Value = "class"
$ws2.Range("E3").Value = "order"
$ws2.Range("F3").Value = "family"
$ws2.Range("G3").Value = "genus"
$ws2.Range("H3").Value = "species"
$ws2.Range("I3").Value = "NA"

# --- Row 4: Brazil -----------------------------------------------------
$ws2.Range("A4").Value = "Brazil"
$ws2.Range("B4").Value = "kingdom"
$ws2.Range("C4").Value = "phylum"
$ws2.Range("D4").Value = "class"
$ws2.Range("E4").Value = "order"
$ws2.Range("F4").Value = "family"
$ws2.Range("G4").Value = "genus"
$ws2.Range("H4").Value = "species, species_guid"
$ws2.Range("I4").Value = "subspecies_guid, subspecies_name"

# --- Row 5: Canada (no data) ---------------------------------------------
$ws2.Range("A5").Value = "Canada"
$ws2.Range("B5:I5").Value = ""

# --- Row 6: Estonia (no data) --------------------------------------------
$ws2.Range("A6").Value = "Estonia"
$ws2.Range("B6:I6").Value = ""

# --- Row 7: France -------------------------------------------------------
$ws2.Range("A7").Value = "France"
$ws2.Range("B7").Value = "kingdom, kingdomID"
$ws2.Range("C7").Value = "phylum, phylumID"
$ws2.Range("D7").Value = "class"
$ws2.Range("E7").Value = "order"
$ws2.Range("F7").Value = "family"
$ws2.Range("G7").Value = "genus"
$ws2.Range("H7").Value = "species, speciesID"
$ws2.Range("I7").Value = "subspecies, subspeciesID"

# --- Row 8: GBIF (not applicable across the board) ------------------------
$ws2.Range("A8").Value = "GBIF"
$ws2.Range("B8:I8").Value = "N/A"

# --- Row 9: Guatemala (no data) ------------------------------------------
$ws2.Range("A9").Value = "Guatemala"
$ws2.Range("B9:I9").Value = ""

# --- Row 10: Portugal (no data) ------------------------------------------
$ws2.Range("A10").Value = "Portugal"
$ws2.Range("B10:I10").Value = ""

# --- Row 11: Spain ---------------------------------------------------------
$ws2.Range("A11").Value = "Spain"
$ws2.Range("B11").Value = "kingdom, kingdom_id"
$ws2.Range("C11").Value = "phylum, phylum_id"
$ws2.Range("D11").Value = "class"
$ws2.Range("E11").Value = "order"
$ws2.Range("F11").Value = "family"
$ws2.Range("G11").Value = "genus"
$ws2.Range("H11").Value = "species, species_guid"
$ws2.Range("I11").Value = "subspecies,subspecies_guid, subspecies_name"

# --- Row 12: Sweden (no data) ---------------------------------------------
$ws2.Range("A12").Value = "Sweden"
$ws2.Range("B12:I12").Value = ""

# --- Row 13: United Kingdom (no data) -------------------------------------
$ws2.Range("A13").Value = "United Kingdom"
$ws2.Range("B13:I13").Value = ""

# ---------------------------------------------------------------------------
# 3) Selections: sheet2 ends on H12; sheet1 stays the active tab, selected
#    over its country list A1:A13.
# ---------------------------------------------------------------------------
$ws2.Range("H12").Select()

$ws1.Activate()
$ws1.Range("A1:A13").Select()
